# Update the "想去人数" (want-to-go count) column (F) values on the
# "展览" and "全部类型" worksheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 62
$wsExpo.Range("F6").Value = 9696
$wsExpo.Range("F7").Value = 873
$wsExpo.Range("F10").Value = 2800
$wsExpo.Range("F13").Value = 20
$wsExpo.Range("F14").Value = 28
$wsExpo.Range("F16").Value = 498
$wsExpo.Range("F17").Value = 103
$wsExpo.Range("F18").Value = 262
$wsExpo.Range("F19").Value = 1387

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 62
$wsAll.Range("F7").Value = 9696
$wsAll.Range("F8").Value = 873
$wsAll.Range("F11").Value = 2800
$wsAll.Range("F14").Value = 20
$wsAll.Range("F15").Value = 28
$wsAll.Range("F17").Value = 498
$wsAll.Range("F18").Value = 103
$wsAll.Range("F19").Value = 262
$wsAll.Range("F20").Value = 1387
